# _Tests.xlsx update:
#  - "Tests" sheet (B1): header renamed from "Exception" to "ExpectedResult"
#  - "Result" sheet (B1): header renamed from "Exception" to "ExpectedResult"
#  - "Result" sheet becomes the active/selected tab (was "Tests")
#  - Selection on "Tests" moved to B14, selection on "Result" moved to D15

$wb = $excel.ActiveWorkbook

$wsTests  = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

# Rename the "expected status" header on both sheets.
$wsTests.Range("B1").Value  = "ExpectedResult"
$wsResult.Range("B1").Value = "ExpectedResult"

# Update the remembered selection on the (no longer active) Tests sheet.
[void]$wsTests.Range("B14").Select()

# Make Result the active sheet/tab and park the selection on D15.
[void]$wsResult.Activate()
[void]$wsResult.Range("D15").Select()
